# Auto-generated edits applying the Chocobo_Profits price-refresh diff.
# For each affected Leve row, update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) with the refreshed market-board values.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 1209
$ws.Range("I28").Value = 1110.5
$ws.Range("K28").Value = 1110.5
$ws.Range("M28").Value = -625.5
# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 197.66667
$ws.Range("I55").Value = 89.5
$ws.Range("J55").Value = 228.57143
$ws.Range("K55").Value = 89.5
$ws.Range("L55").Value = 228.57143
$ws.Range("M55").Value = 124.5
$ws.Range("N55").Value = -656.57143
# Row 123 (Leve Item ID 34090)
$ws.Range("H123").Value = 42980
$ws.Range("J123").Value = 42980
$ws.Range("L123").Value = 42980
$ws.Range("N123").Value = -52780
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 841.22
$ws.Range("I129").Value = 354.3
$ws.Range("J129").Value = 962.95
$ws.Range("K129").Value = 1062.9
$ws.Range("L129").Value = 2888.85
$ws.Range("M129").Value = 3937.1
$ws.Range("N129").Value = -12888.85
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 57780.184
$ws.Range("I132").Value = 62611.03
$ws.Range("J132").Value = 5446
$ws.Range("K132").Value = 187833.09
$ws.Range("L132").Value = 16338
$ws.Range("M132").Value = -185303.09
$ws.Range("N132").Value = -21398
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 3461.0833
$ws.Range("I137").Value = 2824.7
$ws.Range("J137").Value = 6643
$ws.Range("K137").Value = 8474.099999999999
$ws.Range("L137").Value = 19929
$ws.Range("M137").Value = -5924.099999999999
$ws.Range("N137").Value = -25029
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2245.57
$ws.Range("I138").Value = 1132.3214
$ws.Range("K138").Value = 3396.9642
$ws.Range("M138").Value = 1743.0358

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 10382.53
$ws.Range("I32").Value = 6935.66
$ws.Range("J32").Value = 16472
$ws.Range("K32").Value = 6935.66
$ws.Range("L32").Value = 16472
$ws.Range("M32").Value = -6648.66
$ws.Range("N32").Value = -17046
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3714.9
$ws.Range("I122").Value = 1839.8
$ws.Range("K122").Value = 5519.4
$ws.Range("M122").Value = -3069.4
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3407.7036
$ws.Range("I132").Value = 1469.6666
$ws.Range("J132").Value = 4958.1333
$ws.Range("K132").Value = 4408.9998
$ws.Range("L132").Value = 14874.3999
$ws.Range("M132").Value = -1878.9998
$ws.Range("N132").Value = -19934.3999
# Row 137 (Leve Item ID 43227)
$ws.Range("H137").Value = 53265
$ws.Range("J137").Value = 53265
$ws.Range("L137").Value = 53265
$ws.Range("N137").Value = -63465

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 137 (Leve Item ID 42153)
$ws.Range("H137").Value = 34900
$ws.Range("J137").Value = 34900
$ws.Range("L137").Value = 34900
$ws.Range("N137").Value = -45100

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1545.8667
$ws.Range("J16").Value = 2750
$ws.Range("L16").Value = 2750
$ws.Range("N16").Value = -3324
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2821
$ws.Range("I31").Value = 1274.091
$ws.Range("J31").Value = 7075
$ws.Range("K31").Value = 1274.091
$ws.Range("L31").Value = 7075
$ws.Range("M31").Value = -979.0909999999999
$ws.Range("N31").Value = -7665
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2821
$ws.Range("I34").Value = 1274.091
$ws.Range("J34").Value = 7075
$ws.Range("K34").Value = 1274.091
$ws.Range("L34").Value = 7075
$ws.Range("M34").Value = -1072.091
$ws.Range("N34").Value = -7479
# Row 87 (Leve Item ID 11929)
$ws.Range("H87").Value = 23757.143
$ws.Range("J87").Value = 23757.143
$ws.Range("L87").Value = 23757.143
$ws.Range("N87").Value = -26129.143
# Row 90 (Leve Item ID 11929)
$ws.Range("H90").Value = 23757.143
$ws.Range("J90").Value = 23757.143
$ws.Range("L90").Value = 71271.429
$ws.Range("N90").Value = -83127.429
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 3904.8235
$ws.Range("J99").Value = 5076.8887
$ws.Range("L99").Value = 5076.8887
$ws.Range("N99").Value = -8072.8887
# Row 106 (Leve Item ID 18661)
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1545.8667
$ws.Range("J113").Value = 2750
$ws.Range("L113").Value = 2750
$ws.Range("N113").Value = -7090
# Row 116 (Leve Item ID 26117)
$ws.Range("H116").Value = 80000
$ws.Range("J116").Value = 80000
$ws.Range("L116").Value = 80000
$ws.Range("N116").Value = -89178
# Row 118 (Leve Item ID 26112)
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 3904.8235
$ws.Range("J126").Value = 5076.8887
$ws.Range("L126").Value = 15230.6661
$ws.Range("N126").Value = -20170.6661
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 4172.4863
$ws.Range("I134").Value = 4625.852
$ws.Range("J134").Value = 2948.4
$ws.Range("K134").Value = 13877.556
$ws.Range("L134").Value = 8845.200000000001
$ws.Range("M134").Value = -11342.556
$ws.Range("N134").Value = -13915.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 543.4194
$ws.Range("I113").Value = 561.1177
$ws.Range("J113").Value = 521.9286
$ws.Range("K113").Value = 1683.3531
$ws.Range("L113").Value = 1565.7858
$ws.Range("M113").Value = 486.6469
$ws.Range("N113").Value = -5905.7858

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2304.8708
$ws.Range("I132").Value = 1188.1818
$ws.Range("J132").Value = 5034.5557
$ws.Range("K132").Value = 3564.5454
$ws.Range("L132").Value = 15103.6671
$ws.Range("M132").Value = -1034.5454
$ws.Range("N132").Value = -20163.6671
# Row 136 (Leve Item ID 42218)
$ws.Range("H136").Value = 30821.643
$ws.Range("J136").Value = 30821.643
$ws.Range("L136").Value = 92464.929
$ws.Range("N136").Value = -97564.929
# Row 137 (Leve Item ID 43226)
$ws.Range("H137").Value = 86145
$ws.Range("J137").Value = 86145
$ws.Range("L137").Value = 86145
$ws.Range("N137").Value = -96345

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 4822.7046
$ws.Range("I132").Value = 1911.4762
$ws.Range("J132").Value = 7480.7827
$ws.Range("K132").Value = 5734.4286
$ws.Range("L132").Value = 22442.3481
$ws.Range("M132").Value = -3204.4286
$ws.Range("N132").Value = -27502.3481

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 46 (Leve Item ID 42037)
$ws.Range("H46").Value = 86501.875
$ws.Range("J46").Value = 86501.875
$ws.Range("L46").Value = 86501.875
$ws.Range("N46").Value = -86963.875
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 1024.4117
$ws.Range("I107").Value = 967.4167
$ws.Range("J107").Value = 1161.2
$ws.Range("K107").Value = 2902.2501
$ws.Range("L107").Value = 3483.6
$ws.Range("M107").Value = -982.2501000000002
$ws.Range("N107").Value = -7323.6
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 5710.143
$ws.Range("I122").Value = 4312.9546
$ws.Range("J122").Value = 10833.167
$ws.Range("K122").Value = 12938.8638
$ws.Range("L122").Value = 32499.501
$ws.Range("M122").Value = -10488.8638
$ws.Range("N122").Value = -37399.501
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 6413752
$ws.Range("I132").Value = 3806
$ws.Range("J132").Value = 15876054
$ws.Range("K132").Value = 11418
$ws.Range("L132").Value = 47628162
$ws.Range("M132").Value = -8888
$ws.Range("N132").Value = -47633222
# Row 134 (Leve Item ID 42037)
$ws.Range("H134").Value = 86501.875
$ws.Range("J134").Value = 86501.875
$ws.Range("L134").Value = 259505.625
$ws.Range("N134").Value = -264575.625

